# Add "etat commande TMS" column (Q) with header in Q1 and value "valide" in Q2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "etat commande TMS"
$ws.Range("Q2").Value = "valide"

# Widen the new column to match the authored layout (closest value the
# engine's column-width quantization allows to the authored 18.5703125)
$ws.Columns.Item(17).ColumnWidth = 17.67

# Selection moved to Q1:Q2 with active cell Q1 (matches the authored edit)
$ws.Range("Q1:Q2").Select()
